$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Objetivos: fill in the real objectives text (previously misaligned) ---
$ws.Range("B10").Value = "Propiciar ao aluno o conhecimento dos gêneros por meio dos quais ele deverá agir linguisticamente no espaço acadêmico (Objetivo Geral); 2. Ler e redigir resumos acadêmicos e relatórios de pesquisa experimental, além de reconhecer as características de uma resenha (Objetivo Específico); 3. Dominar técnicas de escrita adequadas aos gêneros acadêmicos (Objetivo Específico)."
$ws.Range("C10").Value = $ws.Range("B10").Value2

# --- Insert a new row above "Programa resumido:" so the teacher name lines up
#     with "Docentes responsáveis:" instead of "Programa resumido:" ---
$ws.Rows(13).Insert()

# New row 13 (B/C): teacher name, styled like the other B/C text cells
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B13").Value = "5840514 - Graziela Zamponi"
$ws.Range("C13").Value = $ws.Range("B13").Value2

# --- Row 14 "Programa resumido:" now gets its real short-syllabus text (was "Semestral") ---
$ws.Range("B14").Value = "O texto escrito da esfera acadêmica. Gêneros acadêmicos."
$ws.Range("C14").Value = $ws.Range("B14").Value2

# --- Row 16 "Programa:" gets the full program text (was the stray "01/01/2012") ---
$ws.Range("B16").Value = "1. O texto escrito na esfera acadêmica`n    Aspectos constitutivos do texto escrito`n    Fatores de legibilidade`n    Coesão`n2. Gêneros acadêmicos`n    Noções de gêneros. Gêneros acadêmicos`n    Resumo e resenha`n    Relatório de pesquisa experimental"
$ws.Range("C16").Value = $ws.Range("B16").Value2

# --- Rows 19-21: each one's B/C content shifts up by one slot ---
$ws.Range("B18:C18").ClearContents()
$ws.Range("B19").Value = "N 1  = Prova= 10,0`nN 2 = 1ª NP + 2ª NP  (ver abaixo)"
$ws.Range("C19").Value = $ws.Range("B19").Value2
$ws.Range("B20").Value = "1ª Nota Parcial - Resumo= 5,0`n2ª Nota Parcial - Relatório=5,0 NOTA FINAL = N1 + N2/ 2"
$ws.Range("C20").Value = $ws.Range("B20").Value2
$ws.Range("B21").Value = "Ao aluno que não alcançar a média 5,0 (cinco) no final do período letivo será dada uma recuperação, por meio de uma prova."
$ws.Range("C21").Value = $ws.Range("B21").Value2

# --- Row 22 "Bibliografia:" gets the real bibliography (was the stray recovery-rule text) ---
$ws.Range("B22").Value = "1. FARACO, Carlos Alberto; TEZZA, Cristóvão. Oficina de texto.  6 ed. Petrópolis: Vozes, 2008.`n2. ILARI, Rodolfo.  Introdução à Semântica: brincando com a gramática.  São Paulo: Contexto, 2001.`n3. ______.  Introdução ao estudo do léxico: brincando com as palavras.  São Paulo: Contexto, 2002.`n4. KLEIMAN, Ângela. Texto e leitor: aspectos cognitivos da leitura. 4.ed.  Campinas: Pontes, 1995.`n5. KOCH, Ingedore Villaça.  A coesão textual.  São Paulo: Contexto, 2001.`n6. LIBERATO, Yara; FULGÊNCIO, Lúcia.   É possível facilitar a leitura: um guia para escrever claro.  São Paulo: Contexto, 2007.`n7. MACHADO, A.R (coord.); LOUSADA, E.; ABREU-TARDELLI, L. S.  Resumo.  São Paulo: Parábola Editorial, 2004.`n8. ______.   Resenha.  São Paulo: Parábola Editorial, 2004.`n9. MARCUSCHI, Luiz Antônio.  Da fala para a escrita: atividades de retextualização.  São Paulo: Cortez, 2000.`n10. SERAFINI, Maria José.    Como escrever textos. 5.ed. São Paulo: Globo, 1992."
$ws.Range("C22").Value = $ws.Range("B22").Value2
